$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 28, shifting the existing rows 28-48 down to 29-49.
# Excel's default row-insert copies formatting down from the row above (row 27),
# which carries the date number format (style index 2) on column D.
$ws.Rows.Item(28).Insert()

# Populate the newly inserted row 28 with the new weekly price record.
$ws.Range("A28").Value = 5
$ws.Range("B28").Value = "Macroferia Regional de Talca"
$ws.Range("C28").Value = "Maule"
$ws.Range("D28").Value = 44803
$ws.Range("E28").Value = 7
$ws.Range("F28").Value = 100112043
$ws.Range("G28").Value = "Pepino dulce"
$ws.Range("H28").Value = "Cultivar IV Región"
$ws.Range("I28").Value = "Primera"
$ws.Range("J28").Value = 350
$ws.Range("K28").Value = 14000
$ws.Range("L28").Value = 14000
$ws.Range("M28").Value = 14000
$ws.Range("N28").Value = "$/bandeja 18 kilos"
$ws.Range("O28").Value = "Provincia de Limarí"
$ws.Range("P28").Value = 778
$ws.Range("Q28").Value = 18
$ws.Range("R28").Value = "Hortaliza"
